$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.993.51'
$ws.Range('E2').Value = '  -1.44%  '
$ws.Range('D3').Value = '3.155.80'
$ws.Range('E3').Value = '  -0.95%  '
$ws.Range('E4').Value = '  +0.03%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '610.39'
$ws.Range('E5').Value = '  +1.90%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '147.11'
$ws.Range('E6').Value = '  -3.86%  '
$ws.Range('E7').Value = '  -0.08%  '
$ws.Range('D8').Value = '3.151.81'
$ws.Range('E8').Value = '  -0.98%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.524'
$ws.Range('E9').Value = '  -0.68%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.151'
$ws.Range('E10').Value = '  -2.41%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.38'
$ws.Range('E11').Value = '  -3.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.472'
$ws.Range('E12').Value = '  -1.31%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000254'
$ws.Range('E13').Value = '  -1.92%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '35.58'
$ws.Range('E14').Value = '  -4.17%  '
$ws.Range('D15').Value = '3.652.29'
$ws.Range('E15').Value = '  -1.82%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.117'
$ws.Range('E16').Value = '  +2.91%  '
$ws.Range('D17').Value = '63.925.45'
$ws.Range('E17').Value = '  -1.57%  '
$ws.Range('D18').Value = '3.137.14'
$ws.Range('E18').Value = '  -1.79%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.90'
$ws.Range('E19').Value = '  -2.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '477.01'
$ws.Range('E20').Value = '  -1.48%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.63'
$ws.Range('E21').Value = '  -1.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.716'
$ws.Range('E22').Value = '  -0.45%  '
$ws.Range('E23').Value = '  +2.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.78'
$ws.Range('E24').Value = '  -1.14%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.42'
$ws.Range('E25').Value = '  -2.34%  '
$ws.Range('E26').Value = '  +0.12%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.82'
$ws.Range('E27').Value = '  -4.12%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.54'
$ws.Range('E28').Value = '  -1.83%  '
$ws.Range('B29').Value = 'NEARProtocol'
$ws.Range('C29').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.13'
$ws.Range('E29').Value = '  +2.23%  '
$ws.Range('B30').Value = 'Hedera'
$ws.Range('C30').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.121'
$ws.Range('E30').Value = '  -5.57%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '2.10'
$ws.Range('E31').Value = '  -7.93%  '
$ws.Range('E32').Value = '  +0.00%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.72'
$ws.Range('E33').Value = '  -0.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.34'
$ws.Range('E34').Value = '  -2.32%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.13'
$ws.Range('E35').Value = '  +1.53%  '
$ws.Range('D36').Value = '0.0₃0787'
$ws.Range('E36').Value = '  +7.86%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '6.01'
$ws.Range('E37').Value = '  -2.17%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '52.75'
$ws.Range('E38').Value = '  -3.23%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '459.79'
$ws.Range('E39').Value = '  -3.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.03'
$ws.Range('E40').Value = '  -8.17%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0397'
$ws.Range('E41').Value = '  -1.97%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.119'
$ws.Range('E42').Value = '  -4.83%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.34'
$ws.Range('E43').Value = '  -2.36%  '
$ws.Range('D44').Value = '2.865.47'
$ws.Range('E44').Value = '  -1.77%  '
$ws.Range('B45').Value = 'Fetch.AI'
$ws.Range('C45').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.31'
$ws.Range('E45').Value = '  -4.95%  '
$ws.Range('B46').Value = 'TheGraph'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.269'
$ws.Range('E46').Value = '  -3.07%  '
$ws.Range('E47').Value = '  +2.11%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '26.51'
$ws.Range('E48').Value = '  -3.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.998'
$ws.Range('E49').Value = '  -0.01%  '
$ws.Range('E50').Value = '  -2.40%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '118.92'
$ws.Range('E51').Value = '  -1.83%  '
